# 3DES frequency sheet - aula11: add a new week (columns V:Z) mirroring the
# existing PROJ/PDMO/RMST header block, the week's dates, and each
# student's P/F attendance mark for the new PROJ day (column V).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# --- Row 1: header labels for the new week (PROJ, PROJ, PROJ, PDMO, RMST) ---
$ws.Range("V1").Value = "PROJ"
$ws.Range("W1").Value = "PROJ"
$ws.Range("X1").Value = "PROJ"
$ws.Range("Y1").Value = "PDMO"
$ws.Range("Z1").Value = "RMST"

# --- Row 2: dates for the new week (Mon 01-Feb-2021 .. Fri 05-Feb-2021) ---
$ws.Range("V2").Value = 44228
$ws.Range("V2").NumberFormat = "d-mmm"
$ws.Range("W2").Value = 44229
$ws.Range("W2").NumberFormat = "d-mmm"
$ws.Range("X2").Value = 44230
$ws.Range("X2").NumberFormat = "d-mmm"
$ws.Range("Y2").Value = 44231
$ws.Range("Y2").NumberFormat = "d-mmm"
$ws.Range("Z2").Value = 44232
$ws.Range("Z2").NumberFormat = "d-mmm"

# --- Rows 3-20: attendance (P/F) for the first PROJ session of the new week ---
$ws.Range("V3").Value = "P"
$ws.Range("V4").Value = "P"
$ws.Range("V5").Value = "P"
$ws.Range("V6").Value = "P"
$ws.Range("V7").Value = "F"
$ws.Range("V8").Value = "F"
$ws.Range("V9").Value = "P"
$ws.Range("V10").Value = "F"
$ws.Range("V11").Value = "P"
$ws.Range("V12").Value = "F"
$ws.Range("V13").Value = "P"
$ws.Range("V14").Value = "P"
$ws.Range("V15").Value = "P"
$ws.Range("V16").Value = "P"
$ws.Range("V17").Value = "P"
$ws.Range("V18").Value = "F"
$ws.Range("V19").Value = "P"
$ws.Range("V20").Value = "F"

# --- Column widths: U (now alongside V:Z) widens slightly to fit the new block ---
$ws.Range("U1:Z20").ColumnWidth = 5.8333333

# --- Selection moves to I11 as left by the editor ---
$ws.Range("I11").Select()
